$wb = $excel.ActiveWorkbook

# Sheets: 1 = "102", 2 = "103", 3 = "104", 4 = "105"
$ws104 = $wb.Worksheets.Item(3)
$ws105 = $wb.Worksheets.Item(4)

# --- Sheet "104" (IA-04) : fill in DEVICE_ID / DEVICE_NAME columns ---
$ws104.Range("E7").Value = "BO-07U:PS-CH, BO-07U:PS-CV, BO-09U:PS-CH, BO-09U:PS-CV"

$ws104.Range("D8").Value = "1,2,3,4,5,6,7"
$ws104.Range("E8").Value = "SI-04M2:PS-QFB, SI-04M2:PS-QDB1, SI-04M2:PS-QDB2, SI-05M1:PS-QFA, SI-05M1:PS-QDA, SI-04M1:PS-QS, SI-04M2:PS-QS"

$ws104.Range("D9").Value = "1,2,3,4,5,6,7"
$ws104.Range("E9").Value = "SI-04M1:PS-CH, SI-04M1:PS-CV, SI-04M2:PS-CH, SI-04M2:PS-CV, SI-04C2:PS-CH, SI-04C2:PS-CV-1, SI-04C2:PS-CV-2"

$ws104.Range("D10").Value = "1,2,3,4,5,6,7,8,9,10,11"
$ws104.Range("E10").Value = "SI-04C1:PS-Q1, SI-04C1:PS-Q2, SI-04C2:PS-Q3, SI-04C2:PS-Q4, SI-04C4:PS-Q1, SI-04C4:PS-Q2, SI-04C3:PS-Q3, SI-04C3:PS-Q4, SI-04C1:PS-QS, SI-04C2:PS-QS, SI-04C3:PS-QS"

$ws104.Range("D11").Value = "1,2,3,4,5,6,7"
$ws104.Range("E11").Value = "SI-04C1:PS-CH, SI-04C1:PS-CV, SI-04C4:PS-CH, SI-04C4:PS-CV, SI-04C3:PS-CH, SI-04C3:PS-CV-1, SI-04C3:PS-CV-2"

# --- Sheet "105" (IA-05) : fill in DEVICE_ID / DEVICE_NAME columns ---
$ws105.Range("E9").Value = "BO-11U:PS-CH, BO-11U:PS-CV, BO-13U:PS-CH, BO-13U:PS-CV"

$ws105.Range("D10").Value = "1,2,3,4,5,6,7"
$ws105.Range("E10").Value = "SI-05M2:PS-QFA, SI-05M2:PS-QDA, SI-06M1:PS-QFB, SI-06M1:PS-QDB1, SI-06M1:PS-QDB2, SI-05M1:PS-QS, SI-05M2:PS-QS"

$ws105.Range("D11").Value = "1,2,3,4,5,6,7"
$ws105.Range("E11").Value = "SI-05M1:PS-CH, SI-05M1:PS-CV, SI-05M2:PS-CH, SI-05M2:PS-CV, SI-05C2:PS-CH, SI-05C2:PS-CV-1, SI-05C2:PS-CV-2"

$ws105.Range("D12").Value = "1,2,3,4,5,6,7,8,9,10,11"
$ws105.Range("E12").Value = "SI-05C1:PS-Q1, SI-05C1:PS-Q2, SI-05C2:PS-Q3, SI-05C2:PS-Q4, SI-05C4:PS-Q1, SI-05C4:PS-Q2, SI-05C3:PS-Q3, SI-05C3:PS-Q4, SI-05C1:PS-QS, SI-05C2:PS-QS, SI-05C3:PS-QS"

$ws105.Range("D13").Value = "1,2,3,4,5,6,7"
$ws105.Range("E13").Value = "SI-05C1:PS-CH, SI-05C1:PS-CV, SI-05C4:PS-CH, SI-05C4:PS-CV, SI-05C3:PS-CH, SI-05C3:PS-CV-1, SI-05C3:PS-CV-2"

# --- Update view / selection state on each sheet ---
$ws103 = $wb.Worksheets.Item(2)
$ws103.Range("E8").Select()

$ws105.Range("D12").Select()

$ws104.Activate()
$ws104.Range("C17").Select()
